$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: E1 becomes "costo", F1 becomes "fecha_registro"
$ws.Range("E1").Value = "costo"
$ws.Range("F1").Value = "fecha_registro"

# Update row 2 values
$ws.Range("C2").Value = 1000
$ws.Range("E2").Value = 78000
$ws.Range("F2").Value = "31/1/2026"

# Delete rows 3 and 4 (the old row 3 data and the Cerveza Corona row)
$ws.Range("A3:F4").Delete()
